$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.157479047775269
$ws.Range("B1").Value = 2.367790222167969
$ws.Range("D1").Value = 2.401871204376221
$ws.Range("E1").Value = 1.223013997077942
